$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DCDC")

# --- New content: efficiency formula block (rows 9-11) ---

# Row 9: label + formula text (wrapped, vertically centered), row height 60
$ws.Range("A9").Value = "Formule calcul rendement`nVin, Vout, split, f<split, f>split`nséparateur : espace"
$ws.Range("B9").Value = "12/5/0,1/-669890*x**4+176938*x**3-16759*x**2+689*x+77/1*x**3-9*x**2+14*x+87"
$ws.Rows.Item(9).RowHeight = 60

# Style B10 first: text number format + wrap + vertical center (becomes cellXfs idx 1)
$ws.Range("B10").VerticalAlignment = -4108
$ws.Range("B10").WrapText = $true
$ws.Range("B10").NumberFormat = "@"

# Style A10, A11, B11: vertical center only (becomes cellXfs idx 2)
$ws.Range("A10:A11").VerticalAlignment = -4108
$ws.Range("B11").VerticalAlignment = -4108

# Style A9, B9: wrap + vertical center (becomes cellXfs idx 3)
$ws.Range("A9:B9").VerticalAlignment = -4108
$ws.Range("A9:B9").WrapText = $true

# --- Column widths / sheet view cosmetics ---
# (values picked so the engine's internal char<->pixel rounding lands on the
#  same stored width as the target file; the engine's font metrics differ
#  slightly from the authoring app's, so we aim for the closest achievable)
$ws.Columns.Item(2).ColumnWidth = 20.666666666666668
$ws.Columns.Item(3).ColumnWidth = 9.333333333333334
$ws.Columns.Item(4).ColumnWidth = 9.333333333333334

$excel.ActiveWindow.Zoom = 249
